# Apply cryptos.xlsx data refresh (Fri May 12 05:29:18 UTC 2023 GitHub Actions run)
# Updates Price (col D) and Volume(1h) (col E) values for rows 2-51, and
# restores row order for two coin pairs that were re-ranked (rows 28/29 and
# rows 33/34 swap their Coin/Link/Price/Volume contents while keeping the
# row's rank index in column A untouched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($ws, $addr, $val)
    $cell = $ws.Range($addr)
    # Force text storage so numeric-looking strings (e.g. "1.006") are not
    # reinterpreted as numbers, matching the source data's inline-string type.
    $cell.NumberFormat = "@"
    $cell.Value = $val
    # Drop back to the default style so no stray formatting/quote-prefix
    # marker is introduced (the source cells carry no explicit style).
    $cell.Style = "Normal"
}

Set-CellText $ws 'D2' '26.750.49'
Set-CellText $ws 'E2' '  -2.79%  '
Set-CellText $ws 'D3' '1.775.03'
Set-CellText $ws 'E3' '  -3.11%  '
Set-CellText $ws 'D4' '1.006'
Set-CellText $ws 'E4' '  +0.40%  '
Set-CellText $ws 'D5' '1.005'
Set-CellText $ws 'E5' '  +0.40%  '
Set-CellText $ws 'D6' '306.35'
Set-CellText $ws 'E6' '  -2.04%  '
Set-CellText $ws 'D7' '0.4397'
Set-CellText $ws 'E7' '  +2.56%  '
Set-CellText $ws 'D8' '0.3632'
Set-CellText $ws 'E8' '  -0.81%  '
Set-CellText $ws 'D9' '0.07188'
Set-CellText $ws 'E9' '  -1.30%  '
Set-CellText $ws 'D10' '0.8386'
Set-CellText $ws 'E10' '  -3.03%  '
Set-CellText $ws 'D11' '20.25'
Set-CellText $ws 'E11' '  -1.94%  '
Set-CellText $ws 'D12' '1.824.06'
Set-CellText $ws 'E12' '  -2.40%  '
Set-CellText $ws 'D13' '5.251'
Set-CellText $ws 'E13' '  -2.78%  '
Set-CellText $ws 'D14' '6.365'
Set-CellText $ws 'E14' '  -2.74%  '
Set-CellText $ws 'D15' '0.06827'
Set-CellText $ws 'E15' '  -1.83%  '
Set-CellText $ws 'E16' '  +0.58%  '
Set-CellText $ws 'D17' '79.46'
Set-CellText $ws 'E17' '  -1.55%  '
Set-CellText $ws 'D18' '0.000008694'
Set-CellText $ws 'E18' '  -2.51%  '
Set-CellText $ws 'D19' '1.005'
Set-CellText $ws 'E19' '  +0.39%  '
Set-CellText $ws 'D20' '14.98'
Set-CellText $ws 'E20' '  -2.86%  '
Set-CellText $ws 'D21' '26.669.00'
Set-CellText $ws 'E21' '  -3.48%  '
Set-CellText $ws 'D22' '5.013'
Set-CellText $ws 'E22' '  -2.55%  '
Set-CellText $ws 'E23' '  +1.67%  '
Set-CellText $ws 'D24' '1.974.19'
Set-CellText $ws 'E24' '  -5.93%  '
Set-CellText $ws 'D25' '1.907'
Set-CellText $ws 'D26' '153.53'
Set-CellText $ws 'E26' '  -0.66%  '
Set-CellText $ws 'D27' '18.17'
Set-CellText $ws 'E27' '  -3.88%  '
Set-CellText $ws 'B28' 'InternetComputer(DFINITY)'
Set-CellText $ws 'C28' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-CellText $ws 'D28' '5.045'
Set-CellText $ws 'E28' '  -1.66%  '
Set-CellText $ws 'B29' 'BitcoinCash'
Set-CellText $ws 'C29' 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-CellText $ws 'D29' '114.79'
Set-CellText $ws 'E29' '  +0.45%  '
Set-CellText $ws 'D30' '1.638'
Set-CellText $ws 'E30' '  -10.91%  '
Set-CellText $ws 'D31' '0.08992'
Set-CellText $ws 'E31' '  +1.43%  '
Set-CellText $ws 'D32' '0.7198'
Set-CellText $ws 'E32' '  -4.23%  '
Set-CellText $ws 'B33' 'Filecoin'
Set-CellText $ws 'C33' 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
Set-CellText $ws 'D33' '4.326'
Set-CellText $ws 'E33' '  -4.82%  '
Set-CellText $ws 'B34' 'HuobiToken'
Set-CellText $ws 'C34' 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-CellText $ws 'D34' '2.791'
Set-CellText $ws 'E34' '  -7.11%  '
Set-CellText $ws 'D35' '1.089'
Set-CellText $ws 'E35' '  -3.94%  '
Set-CellText $ws 'D36' '1.005'
Set-CellText $ws 'E36' '  +0.43%  '
Set-CellText $ws 'D37' '1.073'
Set-CellText $ws 'E37' '  -1.36%  '
Set-CellText $ws 'D38' '0.05094'
Set-CellText $ws 'E38' '  -4.40%  '
Set-CellText $ws 'E39' '  -2.68%  '
Set-CellText $ws 'D40' '0.4925'
Set-CellText $ws 'E40' '  -3.13%  '
Set-CellText $ws 'E41' '  -3.22%  '
Set-CellText $ws 'D42' '2.592'
Set-CellText $ws 'E42' '  -7.24%  '
Set-CellText $ws 'D43' '6.128'
Set-CellText $ws 'E43' '  -6.60%  '
Set-CellText $ws 'D44' '7.923'
Set-CellText $ws 'E44' '  -4.94%  '
Set-CellText $ws 'D45' '104.79'
Set-CellText $ws 'E45' '  -1.07%  '
Set-CellText $ws 'D46' '1.006'
Set-CellText $ws 'D47' '10.04'
Set-CellText $ws 'E47' '  -3.85%  '
Set-CellText $ws 'D48' '0.06225'
Set-CellText $ws 'E48' '  -4.18%  '
Set-CellText $ws 'D49' '0.4486'
Set-CellText $ws 'E49' '  -4.26%  '
Set-CellText $ws 'E50' '  -2.21%  '
Set-CellText $ws 'D51' '1.712'
Set-CellText $ws 'E51' '  -0.82%  '
